$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: drop the "2nd" / "3rd" execution columns (H1:M1) ---
$ws.Range("H1:M1").ClearContents()

# --- Row 2: latest execution re-run overwrote the 1st execution result ---
$ws.Range("E2").Value = "18-Nov-2025 02:19:03 PM"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "Pass"

# drop the old 2nd / 3rd execution columns for row 2 as well
$ws.Range("H2:M2").ClearContents()

# --- Rows 3-8: record the latest (re-)execution results ---
$rows = 3..8
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "18-Nov-2025 02:19:03 PM"
    $ws.Range("F$r").Value = ""
    $ws.Range("G$r").Value = "Pass"
}

# --- selection marker moves to G11 (below the used range) ---
[void]$ws.Range("G11").Select()
